$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.337.38'
$ws.Range("E2").Value = '  +0.03%  '

# Row 3
$ws.Range("D3").Value = '1.870.73'
$ws.Range("E3").Value = '  +0.53%  '

# Row 4
$ws.Range("E4").Value = '  +0.04%  '

# Row 5
$ws.Range("D5").Value = '''0.7248'
$ws.Range("E5").Value = '  +2.84%  '

# Row 6
$ws.Range("D6").Value = '''240.95'
$ws.Range("E6").Value = '  +1.17%  '

# Row 7
$ws.Range("E7").Value = '  +0.03%  '

# Row 8
$ws.Range("D8").Value = '''0.07914'
$ws.Range("E8").Value = '  +0.46%  '

# Row 9
$ws.Range("E9").Value = '  +1.33%  '

# Row 10
$ws.Range("D10").Value = '''25.35'
$ws.Range("E10").Value = '  +1.62%  '

# Row 11
$ws.Range("D11").Value = '''0.08258'
$ws.Range("E11").Value = '  +0.93%  '

# Row 12
$ws.Range("B12").Value = 'Polygon'
$ws.Range("C12").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D12").Value = '''0.7236'
$ws.Range("E12").Value = '  +0.68%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.859.01'
$ws.Range("E13").Value = '  -7.98%  '

# Row 14
$ws.Range("D14").Value = '''5.245'
$ws.Range("E14").Value = '  +0.34%  '

# Row 15
$ws.Range("D15").Value = '''90.81'
$ws.Range("E15").Value = '  +1.30%  '

# Row 16
$ws.Range("D16").Value = '29.343.75'
$ws.Range("E16").Value = '  +0.05%  '

# Row 17
$ws.Range("D17").Value = '''5.848'
$ws.Range("E17").Value = '  +0.01%  '

# Row 18
$ws.Range("D18").Value = '''243.71'
$ws.Range("E18").Value = '  +2.13%  '

# Row 19
$ws.Range("E19").Value = '  +0.30%  '

# Row 20
$ws.Range("E20").Value = '  -0.29%  '

# Row 21
$ws.Range("D21").Value = '2.117.97'
$ws.Range("E21").Value = '  -7.66%  '

# Row 22
$ws.Range("E22").Value = '  +0.05%  '

# Row 23
$ws.Range("D23").Value = '''7.987'
$ws.Range("E23").Value = '  +5.15%  '

# Row 24
$ws.Range("D24").Value = '''1.0000'

# Row 25
$ws.Range("D25").Value = '''0.1611'
$ws.Range("E25").Value = '  +12.38%  '

# Row 26
$ws.Range("D26").Value = '''162.27'
$ws.Range("E26").Value = '  -0.34%  '

# Row 27
$ws.Range("D27").Value = '''8.977'
$ws.Range("E27").Value = '  +0.63%  '

# Row 28
$ws.Range("D28").Value = '''18.26'
$ws.Range("E28").Value = '  +0.87%  '

# Row 29
$ws.Range("D29").Value = '''1.350'
$ws.Range("E29").Value = '  -1.87%  '

# Row 30
$ws.Range("E30").Value = '  +1.21%  '

# Row 31
$ws.Range("D31").Value = '''4.376'

# Row 32
$ws.Range("D32").Value = '''4.106'
$ws.Range("E32").Value = '  +1.18%  '

# Row 33
$ws.Range("D33").Value = '''0.05203'
$ws.Range("E33").Value = '  +0.03%  '

# Row 34
$ws.Range("D34").Value = '''1.946'
$ws.Range("E34").Value = '  +2.11%  '

# Row 35
$ws.Range("D35").Value = '''1.187'
$ws.Range("E35").Value = '  +0.57%  '

# Row 36
$ws.Range("E36").Value = '  +1.33%  '

# Row 37
$ws.Range("D37").Value = '''2.677'
$ws.Range("E37").Value = '  +0.04%  '

# Row 38
$ws.Range("D38").Value = '''0.01858'
$ws.Range("E38").Value = '  +0.02%  '

# Row 39
$ws.Range("E39").Value = '  +0.37%  '

# Row 40
$ws.Range("D40").Value = '1.173.95'
$ws.Range("E40").Value = '  -0.22%  '

# Row 41
$ws.Range("D41").Value = '''0.9037'
$ws.Range("E41").Value = '  -2.02%  '

# Row 42
$ws.Range("D42").Value = '''6.127'
$ws.Range("E42").Value = '  +1.65%  '

# Row 43
$ws.Range("D43").Value = '''72.64'
$ws.Range("E43").Value = '  +1.02%  '

# Row 44
$ws.Range("D44").Value = '''1.000'
$ws.Range("E44").Value = '  +0.05%  '

# Row 45
$ws.Range("D45").Value = '''101.78'
$ws.Range("E45").Value = '  -0.43%  '

# Row 46
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").Value = '''0.5281'
$ws.Range("E46").Value = '  -0.86%  '

# Row 47
$ws.Range("B47").Value = 'RocketPoolETH'
$ws.Range("C47").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D47").Value = '2.011.33'
$ws.Range("E47").Value = '  -7.61%  '

# Row 48
$ws.Range("D48").Value = '''1.787'
$ws.Range("E48").Value = '  +1.16%  '

# Row 49
$ws.Range("D49").Value = '''2.899'
$ws.Range("E49").Value = '  +5.86%  '

# Row 50
$ws.Range("D50").Value = '''9.258'
$ws.Range("E50").Value = '  +0.68%  '

# Row 51
$ws.Range("D51").Value = '''0.4285'
$ws.Range("E51").Value = '  +0.17%  '
